$d = $word.ActiveDocument

# --- Change 1: merge the run containing the Greek "alpha" with the run containing "." into a single run ---
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$found1 = $rng1.Find.Execute("desejado " + [char]0x03B1 + ".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find target text for Change 1 (alpha + period run merge)"
}
$mergeTarget = $d.Range($rng1.End - 2, $rng1.End)
$mergeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>α.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$mergeTarget.InsertXML($mergeXml, "Replace")

# --- Change 2a: insert new explanatory-text runs right before the _GoBack bookmark (same paragraph) ---
$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Range.Start
$insertRunsTarget = $d.Range($bmStart, $bmStart)
$runsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Agora existe um termo em que você deve se familiarizar, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:b/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Margem de Erro</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">. É dada como </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>{(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>z.σ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>)/</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>√n}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> e é definida exemplificando o erro pela supervisão da pessoa que coletou os exemplos. Isso significa, se uma média amostral se encontrar no nível da margem de erro </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertRunsTarget.InsertXML($runsXml)

# --- Change 2b: insert two new empty paragraphs right after the (now-extended) bookmark paragraph ---
$bm2 = $d.Bookmarks("_GoBack")
$bmStart2 = $bm2.Range.Start
$insertParasTarget = $d.Range($bmStart2, $bmStart2)
$twoParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="792"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="792"/><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="36394D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertParasTarget.InsertXML($twoParasXml)

Write-Host "Done"
